$d = $word.ActiveDocument

# Update the date heading (first paragraph)
$d.Content.Find.Execute("2024-08-17 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-08-18 Sunday", 2)

# Update the division problems in the table, cell by cell (row, col are 1-indexed)
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "103÷4=25, 3"
$t.Cell(1,2).Range.Text  = "845÷7=120, 5"
$t.Cell(1,3).Range.Text  = "744÷3=248, 0"
$t.Cell(1,4).Range.Text  = "394÷8=49, 2"
$t.Cell(1,5).Range.Text  = "728÷3=242, 2"

$t.Cell(5,1).Range.Text  = "515÷7=73, 4"
$t.Cell(5,2).Range.Text  = "750÷6=125, 0"
$t.Cell(5,3).Range.Text  = "906÷6=151, 0"
$t.Cell(5,4).Range.Text  = "879÷2=439, 1"
$t.Cell(5,5).Range.Text  = "583÷2=291, 1"

$t.Cell(9,1).Range.Text  = "897÷3=299, 0"
$t.Cell(9,2).Range.Text  = "958÷8=119, 6"
$t.Cell(9,3).Range.Text  = "167÷9=18, 5"
$t.Cell(9,4).Range.Text  = "112÷4=28, 0"
$t.Cell(9,5).Range.Text  = "128÷3=42, 2"

$t.Cell(13,1).Range.Text = "882÷2=441, 0"
$t.Cell(13,2).Range.Text = "680÷4=170, 0"
$t.Cell(13,3).Range.Text = "281÷2=140, 1"
$t.Cell(13,4).Range.Text = "600÷4=150, 0"
$t.Cell(13,5).Range.Text = "653÷5=130, 3"

$t.Cell(17,1).Range.Text = "609÷2=304, 1"
$t.Cell(17,2).Range.Text = "245÷2=122, 1"
$t.Cell(17,3).Range.Text = "711÷6=118, 3"
$t.Cell(17,4).Range.Text = "722÷6=120, 2"
$t.Cell(17,5).Range.Text = "783÷8=97, 7"
